$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header A1: "Floor" -> "Description"
$ws.Range("A1").Value = "Description"

# Row 2, column A: empty -> "changed from Python Script"
$ws.Range("A2").Value = "changed from Python Script"

# Update the selection to A3 (matches <selection activeCell="A3" sqref="A3"/>)
$ws.Range("A3").Select()
